# Generate Report for Handoff
#
# A new handoff batch was generated for the 2f1ea5bb / 50120022 /
# a0e67254 / e608a588 files (rows 4-7 of the per-locale tables), so the
# localization-status report is regenerated for that batch:
#
#  - zh-cn: those four files are now handed off, so Priority flips from
#    "low" to "ht" and the Latest Handoff Datetime is refreshed
#    (rows 4-7, cols E & H).
#  - de-de: same four files flip Priority from "low" to "ht" and get a
#    refreshed Latest Handoff Datetime too (rows 4-7, cols E & H).
#  - Overview: its "Latest HO Xliff Generate Date" column mirrors the
#    de-de handoff timestamp, so it picks up the same refreshed value
#    (rows 4-7, col G).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-13 06:37:44"

    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-13 06:37:52"

    $wsOverview.Cells.Item($r, 7).Value = "2016-08-13 06:37:52"
}
